$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the mismatched row pairs (data was shifted between adjacent match rows).
$rA = $ws.Range("B35:AC35")
$rB = $ws.Range("B36:AC36")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value = $vB
$rB.Value = $vA

$rA = $ws.Range("B38:AC38")
$rB = $ws.Range("B39:AC39")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value = $vB
$rB.Value = $vA

$rA = $ws.Range("B65:AC65")
$rB = $ws.Range("B66:AC66")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value = $vB
$rB.Value = $vA

$rA = $ws.Range("B71:AC71")
$rB = $ws.Range("B72:AC72")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value = $vB
$rB.Value = $vA

$rA = $ws.Range("B125:AC125")
$rB = $ws.Range("B126:AC126")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value = $vB
$rB.Value = $vA

$rA = $ws.Range("B128:AC128")
$rB = $ws.Range("B129:AC129")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value = $vB
$rB.Value = $vA

$rA = $ws.Range("B149:AC149")
$rB = $ws.Range("B150:AC150")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value = $vB
$rB.Value = $vA

$rA = $ws.Range("B167:AC167")
$rB = $ws.Range("B168:AC168")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value = $vB
$rB.Value = $vA

$rA = $ws.Range("B187:AC187")
$rB = $ws.Range("B188:AC188")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value = $vB
$rB.Value = $vA

$rA = $ws.Range("B212:AC212")
$rB = $ws.Range("B213:AC213")
$vA = $rA.Value2
$vB = $rB.Value2
$rA.Value = $vB
$rB.Value = $vA

# Remove the trailing 5 rows that were dropped from the data set.
$ws.Range("A240:A244").EntireRow.Delete()
